$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final roster (player, positions, team) in the new row order. The sheet
# already holds 17 player rows (rows 2-18); the updated roster has 18
# players (row 19 is new) - "Jerami Grant" was dropped and "Eric Gordon" /
# "Grayson Allen" were added, and the whole table was re-sorted into this
# order.
$roster = @(
  @("Dennis Schröder", "PG,SG", "Golden State Warriors"),
  @("Jordan Poole", "PG,SG", "Washington Wizards"),
  @("Jamal Murray", "PG,SG", "Denver Nuggets"),
  @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
  @("Tobias Harris", "SF,PF", "Detroit Pistons"),
  @("RJ Barrett", "SG,SF,PF", "Toronto Raptors"),
  @("Eric Gordon", "SG,SF", "Philadelphia 76ers"),
  @("Grayson Allen", "PG,SG,SF", "Phoenix Suns"),
  @("Bam Adebayo", "C", "Miami Heat"),
  @("Nikola Jovic", "PF,C", "Miami Heat"),
  @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
  @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
  @("Toumani Camara", "SF,PF", "Portland Trail Blazers"),
  @("Mike Conley", "PG", "Minnesota Timberwolves"),
  @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
  @("Zach LaVine", "SG,SF", "Chicago Bulls"),
  @("Joel Embiid", "C", "Philadelphia 76ers"),
  @("John Collins", "PF,C", "Utah Jazz")
)

$row = 2
foreach ($player in $roster) {
  $ws.Cells.Item($row, 1).Value = $player[0]
  $ws.Cells.Item($row, 2).Value = $player[1]
  $ws.Cells.Item($row, 3).Value = $player[2]
  $row = $row + 1
}
